$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.939.14'
$ws.Range("E2").Value = '  -2.01%  '
$ws.Range("D3").Value = '2.726.80'
$ws.Range("E3").Value = '  -6.12%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '''508.02'
$ws.Range("E5").Value = '  -3.59%  '
$ws.Range("D6").Value = '''141.49'
$ws.Range("E6").Value = '  -0.44%  '
$ws.Range("E7").Value = '  -0.35%  '
$ws.Range("D8").Value = '''0.533'
$ws.Range("E8").Value = '  -3.45%  '
$ws.Range("D9").Value = '2.736.72'
$ws.Range("E9").Value = '  -5.82%  '
$ws.Range("D10").Value = '''6.12'
$ws.Range("E10").Value = '  +4.41%  '
$ws.Range("E11").Value = '  -2.09%  '
$ws.Range("D12").Value = '''0.349'
$ws.Range("E12").Value = '  -1.35%  '
$ws.Range("E13").Value = '  +1.49%  '
$ws.Range("D14").Value = '3.200.48'
$ws.Range("E14").Value = '  -6.15%  '
$ws.Range("D15").Value = '58.823.75'
$ws.Range("E15").Value = '  -2.55%  '
$ws.Range("D16").Value = '''21.75'
$ws.Range("E16").Value = '  -3.85%  '
$ws.Range("D17").Value = '''0.0000136'
$ws.Range("E17").Value = '  -2.28%  '
$ws.Range("D18").Value = '2.721.63'
$ws.Range("E18").Value = '  -5.91%  '
$ws.Range("D19").Value = '''4.74'
$ws.Range("E19").Value = '  -3.63%  '
$ws.Range("D20").Value = '''10.97'
$ws.Range("E20").Value = '  -4.33%  '
$ws.Range("D21").Value = '''347.53'
$ws.Range("E21").Value = '  -3.42%  '
$ws.Range("D22").Value = '''6.25'
$ws.Range("E22").Value = '  -4.87%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").Value = '''5.61'
$ws.Range("E24").Value = '  -0.43%  '
$ws.Range("D25").Value = '''62.83'
$ws.Range("E25").Value = '  -0.62%  '
$ws.Range("D26").Value = '''0.425'
$ws.Range("E26").Value = '  -5.12%  '
$ws.Range("D27").Value = '''0.171'
$ws.Range("E27").Value = '  -0.71%  '
$ws.Range("D28").Value = '''0.992'
$ws.Range("E28").Value = '  -0.79%  '
$ws.Range("D29").Value = '0.0₃0840'
$ws.Range("E29").Value = '  -0.82%  '
$ws.Range("D30").Value = '''7.48'
$ws.Range("E30").Value = '  -3.10%  '
$ws.Range("E31").Value = '  -0.19%  '
$ws.Range("D32").Value = '''1.61'
$ws.Range("E32").Value = '  -3.41%  '
$ws.Range("D33").Value = '''19.11'
$ws.Range("E33").Value = '  -1.55%  '
$ws.Range("D34").Value = '''148.90'
$ws.Range("E34").Value = '  -1.41%  '
$ws.Range("D35").Value = '''4.21'
$ws.Range("E35").Value = '  -2.05%  '
$ws.Range("D36").Value = '''5.38'
$ws.Range("E36").Value = '  -2.42%  '
$ws.Range("D37").Value = '''0.947'
$ws.Range("E37").Value = '  -3.88%  '
$ws.Range("D38").Value = '''1.14'
$ws.Range("E38").Value = '  -4.03%  '
$ws.Range("D39").Value = '''36.11'
$ws.Range("E39").Value = '  -4.13%  '
$ws.Range("D40").Value = '''1.39'
$ws.Range("E40").Value = '  -3.80%  '
$ws.Range("D41").Value = '2.190.61'
$ws.Range("E41").Value = '  -6.24%  '
$ws.Range("D42").Value = '''3.53'
$ws.Range("E42").Value = '  -3.21%  '
$ws.Range("B43").Value = 'Hedera'
$ws.Range("C43").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D43").Value = '''0.0557'
$ws.Range("E43").Value = '  -1.76%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Value = '''0.995'
$ws.Range("E44").Value = '  -0.47%  '
$ws.Range("D45").Value = '''0.601'
$ws.Range("E45").Value = '  -6.51%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '''19.08'
$ws.Range("E46").Value = '  -7.61%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '''4.82'
$ws.Range("E47").Value = '  +0.36%  '
$ws.Range("D48").Value = '''10.36'
$ws.Range("E48").Value = '  -0.03%  '
$ws.Range("D49").Value = '''0.0229'
$ws.Range("E49").Value = '  -1.40%  '
$ws.Range("D50").Value = '''0.0887'
$ws.Range("E50").Value = '  -4.22%  '
$ws.Range("D51").Value = '''18.03'
$ws.Range("E51").Value = '  -0.87%  '
